$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data update for Q1 advanced estimate: each line is Row,Col,NewValue
$updates = @"
205,8,327773.6666666667
206,2,257.3506666666667
206,6,4988.666666666667
206,7,13686
206,8,324580.6666666667
207,2,255.0148888888889
207,6,4986.777777777777
207,7,13692.66666666667
207,8,323780.8888888889
208,2,255.681962962963
208,6,4970.259259259259
208,7,13678
208,8,325378.4074074074
209,2,256.0158395061728
209,6,4981.901234567901
209,7,13685.55555555556
209,8,324579.987654321
210,2,255.5708971193415
210,6,4979.64609053498
210,7,13685.40740740741
210,8,324579.7613168724
211,2,255.7562331961591
211,6,4977.268861454047
211,7,13682.98765432099
211,8,324846.0521262003
212,2,255.7809899405578
212,6,4979.605395518976
212,7,13684.65020576132
212,8,324668.6003657979
213,2,255.7027067520195
213,6,4978.840115836
213,7,13684.34842249657
213,8,324698.1379362902
214,2,255.7466432962455
214,6,4978.571457603008
214,7,13683.99542752629
214,8,324737.5968094295
215,2,255.7434466629409
215,6,4979.005656319328
215,7,13684.33135192806
215,8,324701.4450371725
216,2,255.7309322370686
216,6,4978.805743252779
216,7,13684.22506731698
216,8,324712.3932609641
217,2,255.740340732085
217,6,4978.794285725038
217,7,13684.18394892378
217,8,324717.1450358554
218,2,255.7382398773649
218,6,4978.868561765715
218,7,13684.2467893896
218,8,324710.3277779973
219,2,255.7365042821728
219,6,4978.822863581177
219,7,13684.21860187678
219,8,324713.2886916056
220,2,255.7383616305409
220,6,4978.82857035731
220,7,13684.21644673005
220,8,324713.5871684861
221,2,255.7377019300262
221,6,4978.839998568067
221,7,13684.22727933215
221,8,324712.4012126963
222,2,255.7375226142466
222,6,4978.830477502185
222,7,13684.22077597966
222,8,324713.092357596
223,2,255.7378620582713
223,6,4978.833015475854
223,7,13684.22150068062
223,8,324713.0269129261
224,2,255.7376955341814
224,6,4978.834497182035
224,7,13684.22318533081
224,8,324712.8401610728
225,2,255.7376934022331
225,6,4978.832663386692
225,7,13684.2218206637
225,8,324712.9864771983
226,2,255.7377503315619
226,6,4978.833392014861
226,7,13684.22216889171
226,8,324712.9511837324
227,2,255.7377130893254
227,6,4978.833517527863
227,7,13684.22239162874
227,8,324712.9259406679
228,2,255.7377189410401
228,6,4978.833190976472
228,7,13684.22212706138
228,8,324712.9545338662
229,2,255.7377274539758
229,6,4978.833366839732
229,7,13684.22222919394
229,8,324712.9438860889
230,2,255.7377198281138
230,6,4978.833358448022
230,7,13684.22224929469
230,8,324712.941453541
231,2,255.7377220743766
231,6,4978.833305421409
231,7,13684.22220185
231,8,324712.9466244987
232,2,255.7377231188221
232,6,4978.833343569721
232,7,13684.22222677954
232,8,324712.9439880428
233,2,255.7377216737708
233,6,4978.833335813051
233,7,13684.22222597475
233,8,324712.9440220275
234,2,255.7377222889898
234,6,4978.833328268061
234,7,13684.22221820143
234,8,324712.9448781897
235,2,255.7377223605276
235,6,4978.833335883611
235,7,13684.22222365191
235,8,324712.9442960867
236,2,255.7377221077627
236,6,4978.833333321574
236,7,13684.22222260936
236,8,324712.944398768
237,2,255.7377222524267
237,6,4978.833332491082
237,7,13684.22222148757
237,8,324712.9445243481
238,2,255.737722240239
238,6,4978.833333898756
238,7,13684.22222258295
238,8,324712.9444064009
239,2,255.7377222001428
239,6,4978.833333237138
239,7,13684.22222222662
239,8,324712.9444431723
240,2,255.7377222309362
240,6,4978.833333208992
240,7,13684.22222209904
240,8,324712.9444579738
241,2,255.7377222237726
241,6,4978.833333448295
241,7,13684.22222230287
241,8,324712.944435849
242,2,255.7377222182839
242,6,4978.833333298141
242,7,13684.22222220951
242,8,324712.944445665
243,2,255.7377222243309
243,6,4978.833333318476
243,7,13684.22222220381
243,8,324712.9444464959
244,2,255.7377222221291
244,6,4978.833333354971
244,7,13684.22222223873
244,8,324712.94444267
245,2,255.7377222215813
245,6,4978.833333323863
245,7,13684.22222221735
245,8,324712.9444449436
246,2,255.7377222226805
246,6,4978.833333332436
246,7,13684.22222221996
246,8,324712.9444447032
247,2,255.7377222221303
247,6,4978.83333333709
247,7,13684.22222222535
247,8,324712.9444441056
248,2,255.7377222221307
248,6,4978.833333331129
248,7,13684.22222222089
248,8,324712.9444445841
249,2,255.7377222223138
249,6,4978.833333333552
249,7,13684.22222222207
249,8,324712.9444444643
"@

$rows = $updates -split "`n"
foreach ($line in $rows) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $v = [double]$parts[2]
    $ws.Cells.Item($r, $c).Value = $v
}
